$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.807.48'
$ws.Range('E2').Value = '  -0.27%  '
$ws.Range('D3').Value = '2.346.85'
$ws.Range('E3').Value = '  -0.47%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '239.35'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.08%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.659'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -4.72%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '72.70'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -5.89%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.595'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -6.09%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.100'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.09%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '58.48'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.90%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '32.79'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -2.63%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '7.25'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -3.99%  '
$ws.Range('B14').Value = 'TRON'
$ws.Range('C14').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.107'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.66%  '
$ws.Range('D15').Value = '2.696.66'
$ws.Range('E15').Value = '  -0.27%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '16.11'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -5.37%  '
$ws.Range('E17').Value = '  -3.71%  '
$ws.Range('D18').Value = '2.341.25'
$ws.Range('E18').Value = '  -0.73%  '
$ws.Range('D19').Value = '43.740.17'
$ws.Range('E19').Value = '  -0.16%  '
$ws.Range('E20').Value = '  -1.64%  '
$ws.Range('E21').Value = '  -0.93%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '77.78'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.01%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '249.06'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -4.16%  '
$ws.Range('E24').Value = '  +0.13%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.72'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.31%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.86'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.60%  '
$ws.Range('E27').Value = '  -2.19%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.37'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -5.98%  '
$ws.Range('E29').Value = '  -1.39%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '176.98'
$ws.Range('D30').Style = 'Normal'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '22.24'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -3.86%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.126'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.84%  '
$ws.Range('E33').Value = '  -2.37%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0750'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.31%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.07'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -5.92%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.33'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.64%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.75'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.29%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.39'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.81%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.36'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -3.30%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.67'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +24.90%  '
$ws.Range('E41').Value = '  -3.33%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '65.96'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +16.65%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '9.23'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.98%  '
$ws.Range('E44').Value = '  -3.42%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '18.85'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.81%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.195'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -6.43%  '
$ws.Range('E47').Value = '  -0.07%  '
$ws.Range('E48').Value = '  -3.53%  '
$ws.Range('E49').Value = '  -5.55%  '
$ws.Range('E50').Value = '  -3.97%  '
$ws.Range('E51').Value = '  +3.04%  '
